# Update gh-pages to output generated at 456a3b4
# Applies the data refresh to 北京-漫展信息.xlsx: "想去人数"/"最低票价" counter bumps
# on sheets 展览 / 演出 / 全部类型, plus two newly scraped 演出 events inserted
# into the 演出 sheet (shifting the later rows down by two).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "展览": bump "想去人数"/"最低票价" counters for several events.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 2972
$ws1.Range("F3").Value = 6418
$ws1.Range("F4").Value = 2531
$ws1.Range("G6").Value = 49
$ws1.Range("F9").Value = 2933
$ws1.Range("F10").Value = 355
$ws1.Range("F12").Value = 7488
$ws1.Range("F13").Value = 349
$ws1.Range("F16").Value = 252
$ws1.Range("F19").Value = 9138
$ws1.Range("F25").Value = 26
$ws1.Range("F27").Value = 115
$ws1.Range("F30").Value = 108
$ws1.Range("F31").Value = 71
$ws1.Range("F32").Value = 115
$ws1.Range("F37").Value = 1485
$ws1.Range("F38").Value = 765
$ws1.Range("F39").Value = 3918
$ws1.Range("F40").Value = 208
$ws1.Range("F42").Value = 1197
$ws1.Range("F43").Value = 88
$ws1.Range("F45").Value = 242
$ws1.Range("F47").Value = 59

# ---------------------------------------------------------------------------
# Sheet "演出": bump two counters, then insert two freshly scraped events
# right after row 14 (2024-05-12 "母亲节专题视听音乐会"), pushing the
# former rows 15-19 down to rows 17-21.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F7").Value = 126
$ws2.Range("F9").Value = 3

# Insert the two new rows at position 15 (row 14 stays put).
$ws2.Rows.Item(15).Insert()
$ws2.Rows.Item(15).Insert()

function Set-EventRow {
    param($Sheet, $Row, $Index, $StartDate, $Name, $Venue, $TimeRange, $WantCount, $MinPrice, $Link, $Cover)

    $Sheet.Range("A$Row").Value = $Index
    $Sheet.Range("A$Row").Font.Bold = $true
    $Sheet.Range("A$Row").HorizontalAlignment = -4108
    $Sheet.Range("A$Row").VerticalAlignment = -4160
    $Sheet.Range("A$Row").Borders.LineStyle = 1

    # Force the B column (e.g. "2024-05-12") to stay plain text instead of
    # being auto-parsed into a date serial, then drop back to the default
    # "Normal" style so no stray NumberFormat/quote-prefix style sticks.
    $Sheet.Range("B$Row").NumberFormat = "@"
    $Sheet.Range("B$Row").Value = $StartDate
    $Sheet.Range("B$Row").Style = "Normal"

    $Sheet.Range("C$Row").Value = $Name
    $Sheet.Range("D$Row").Value = $Venue
    $Sheet.Range("E$Row").Value = $TimeRange
    $Sheet.Range("F$Row").Value = $WantCount
    $Sheet.Range("G$Row").Value = $MinPrice
    $Sheet.Range("H$Row").Value = $Link
    $Sheet.Range("I$Row").Value = $Cover
}

Set-EventRow $ws2 15 14 "2024-05-12" `
    "北京·《家庭教师》《七龙珠》《火影忍者》超燃动漫音乐会" `
    "北京东图剧场 北京东图剧场" `
    "2024.05.12 19:30-05.12 21:00" `
    0 64 `
    "https://show.bilibili.com/platform/detail.html?id=84067" `
    "//i0.hdslb.com/bfs/openplatform/202404/9gLpckTZ1712754110725.jpeg"

Set-EventRow $ws2 16 15 "2024-05-12" `
    "北京·摇滚新星企划" `
    "朝阳北路甲27号菁英梦谷·常营文创产业园南门B5座 WeShow Live 北京" `
    "2024.05.12 15:00-05.12 17:50" `
    0 99 `
    "https://show.bilibili.com/platform/detail.html?id=84069" `
    "//i2.hdslb.com/bfs/openplatform/202404/1Hz7WZo91712620004229.jpeg"

# Renumber the "A" index column for the rows that shifted down (old rows
# 15-19, now at 17-21) so the running counter stays 0-based/contiguous.
for ($r = 17; $r -le 21; $r++) {
    $ws2.Range("A$r").Value = $r - 1
}

# ---------------------------------------------------------------------------
# Sheet "全部类型": same counter bumps as sheet "展览" (combined listing),
# at its own row offsets.
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 2972
$ws4.Range("F6").Value = 6418
$ws4.Range("F7").Value = 2531
$ws4.Range("F8").Value = 126
$ws4.Range("G10").Value = 49
$ws4.Range("F13").Value = 2933
$ws4.Range("F14").Value = 355
$ws4.Range("F18").Value = 7488
$ws4.Range("F19").Value = 349
$ws4.Range("F22").Value = 252
$ws4.Range("F24").Value = 9138
$ws4.Range("F27").Value = 26
$ws4.Range("F29").Value = 115
$ws4.Range("F31").Value = 108
$ws4.Range("F32").Value = 71
$ws4.Range("F33").Value = 115
$ws4.Range("F37").Value = 1485
$ws4.Range("F38").Value = 765
$ws4.Range("F40").Value = 3918
$ws4.Range("F41").Value = 208
$ws4.Range("F44").Value = 1197
$ws4.Range("F45").Value = 88
$ws4.Range("F46").Value = 242
$ws4.Range("F47").Value = 59
